# Adds a new "file/clue" entry (row 10) to the character sheet, and
# separates it from the previous block (row 9) with a thin bottom border.
# Matches the upstream commit "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Close off the previous group (row 9) with a thin bottom border ---
# (A9 is currently blank but still picks up the new bordered style.)
$prevGroup = $ws.Range("A9:E9")
$prevGroup.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# --- New row 10 content ---
# Order of assignment controls shared-string index allocation, so write
# these in the same order the original authoring tool produced them in.
$ws.Range("C10").Value = ' Then we failed to find clues\nhere too…'
$ws.Range("A10").Value = 'SCRIPT/G01P03A/um1104.ssb'
$ws.Range("D10").Value = ' И здесь мы ничего не нашли...'
$ws.Range("E10").Value = ' É èäåòû íú îéœåãï îå îàšìé...'
$ws.Range("B10").Value = 66

# New row gets a thin border framing it top and bottom.
$newGroup = $ws.Range("A10:E10")
$newGroup.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$newGroup.Borders.Item(8).LineStyle = 1    # xlEdgeTop

# Taller row to fit the wrapped text.
$ws.Rows.Item(10).RowHeight = 26.4

# Move the active selection, as in the saved workbook.
$ws.Range("A3").Select() | Out-Null
